$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "End Date" row (row 5) for Senate (B), House (C), and Total (D) columns
# is updated from 2/22/2002 to 11/22/2002 (serial date 37582).
$ws.Range("B5:D5").Value = 37582
